# Refitting NCDEs to individual patients (for manuscript figure)
#
# Adds a new "Label" column (H) to the classification-results sheet (the
# true Control/MDD class per patient) and refreshes the D/E/F prediction
# values that shifted slightly after refitting the NCDEs per-patient.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Label" header in H1, matching the header style used by B1:G1 ---
$ws.Cells.Item(1, 8).Value = "Label"
$ws.Range("G1").Copy() | Out-Null
$ws.Cells.Item(1, 8).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Block 1: Iterations = 100 (rows 2-11) ---
$ws.Cells.Item(2, 4).Value  = 0.2427961003717462
$ws.Cells.Item(2, 5).Value  = 0.2427961003717462
$ws.Cells.Item(2, 8).Value  = 0

$ws.Cells.Item(3, 4).Value  = 0.2475000994943404
$ws.Cells.Item(3, 5).Value  = 0.2475000994943404
$ws.Cells.Item(3, 8).Value  = 0

$ws.Cells.Item(4, 4).Value  = 0.4887590202610154
$ws.Cells.Item(4, 5).Value  = 0.4887590202610154
$ws.Cells.Item(4, 8).Value  = 0

$ws.Cells.Item(5, 4).Value  = 0.5366789860028048
$ws.Cells.Item(5, 5).Value  = 0.5366789860028048
$ws.Cells.Item(5, 8).Value  = 0

$ws.Cells.Item(6, 4).Value  = 0.7655070637160785
$ws.Cells.Item(6, 5).Value  = 0.7655070637160785
$ws.Cells.Item(6, 8).Value  = 0

$ws.Cells.Item(7, 4).Value  = 0.4981118683992572
$ws.Cells.Item(7, 5).Value  = 0.5018881316007429
$ws.Cells.Item(7, 8).Value  = 1

$ws.Cells.Item(8, 4).Value  = 0.5294824286626919
$ws.Cells.Item(8, 5).Value  = 0.4705175713373081
$ws.Cells.Item(8, 8).Value  = 1

$ws.Cells.Item(9, 4).Value  = 0.5259839304318392
$ws.Cells.Item(9, 5).Value  = 0.4740160695681608
$ws.Cells.Item(9, 8).Value  = 1

$ws.Cells.Item(10, 4).Value = 0.6022086133330473
$ws.Cells.Item(10, 5).Value = 0.3977913866669527
$ws.Cells.Item(10, 8).Value = 1

$ws.Cells.Item(11, 4).Value = 0.5151537931624658
$ws.Cells.Item(11, 5).Value = 0.4848462068375342
$ws.Cells.Item(11, 6).Value = 0.6598767638206482
$ws.Cells.Item(11, 8).Value = 1

# --- Block 2: Iterations = 200 (rows 12-21); D/E/F unchanged, only Label added ---
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(17, 8).Value = 1
$ws.Cells.Item(18, 8).Value = 1
$ws.Cells.Item(19, 8).Value = 1
$ws.Cells.Item(20, 8).Value = 1
$ws.Cells.Item(21, 8).Value = 1
